$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as literal text,
# matching the source file where every data cell is an inline/shared string.
function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).NumberFormat = "General"
    $ws.Range($cellRef).Style = "Normal"
}

$ws.Range('D2').Value = '60.765.38'
$ws.Range('E2').Value = '  -2.84%  '

$ws.Range('D3').Value = '3.362.42'
$ws.Range('E3').Value = '  -0.54%  '

$ws.Range('E4').Value = '  +0.09%  '

Set-TextValue 'D5' '570.97'
$ws.Range('E5').Value = '  -0.30%  '

Set-TextValue 'D6' '135.25'
$ws.Range('E6').Value = '  +7.77%  '

$ws.Range('D8').Value = '3.362.14'
$ws.Range('E8').Value = '  -0.54%  '

Set-TextValue 'D9' '0.477'
$ws.Range('E9').Value = '  +0.41%  '

$ws.Range('E10').Value = '  +4.26%  '

$ws.Range('E11').Value = '  +2.54%  '

$ws.Range('E12').Value = '  +3.84%  '

$ws.Range('D13').Value = '3.934.14'
$ws.Range('E13').Value = '  -0.51%  '

$ws.Range('E14').Value = '  +1.66%  '

$ws.Range('E15').Value = '  +0.97%  '

$ws.Range('D16').Value = '3.364.10'
$ws.Range('E16').Value = '  -0.39%  '

Set-TextValue 'D17' '25.15'
$ws.Range('E17').Value = '  +3.08%  '

$ws.Range('D18').Value = '60.962.77'
$ws.Range('E18').Value = '  -2.56%  '

Set-TextValue 'D19' '13.89'
$ws.Range('E19').Value = '  +6.22%  '

Set-TextValue 'D20' '5.80'
$ws.Range('E20').Value = '  +3.21%  '

Set-TextValue 'D21' '9.42'
$ws.Range('E21').Value = '  +1.95%  '

Set-TextValue 'D22' '372.04'

$ws.Range('E23').Value = '  +2.53%  '

$ws.Range('D24').Value = '3.496.20'
$ws.Range('E24').Value = '  -0.52%  '

$ws.Range('E26').Value = '  -1.44%  '

Set-TextValue 'D27' '0.0000116'
$ws.Range('E27').Value = '  +10.36%  '

$ws.Range('E28').Value = '  +22.24%  '

$ws.Range('E29').Value = '  +11.10%  '

$ws.Range('E30').Value = '  +0.33%  '

Set-TextValue 'D31' '8.10'
$ws.Range('E31').Value = '  +4.28%  '

$ws.Range('E32').Value = '  +1.54%  '

Set-TextValue 'D33' '0.154'
$ws.Range('E33').Value = '  +3.95%  '

$ws.Range('E34').Value = '  -0.07%  '

$ws.Range('D35').Value = '3.391.97'
$ws.Range('E35').Value = '  -0.53%  '

$ws.Range('E36').Value = '  +3.13%  '

Set-TextValue 'D37' '5.55'
$ws.Range('E37').Value = '  +5.02%  '

$ws.Range('E38').Value = '  +4.29%  '

$ws.Range('E39').Value = '  +4.94%  '

Set-TextValue 'D40' '163.02'
$ws.Range('E40').Value = '  -2.19%  '

Set-TextValue 'D41' '0.0784'
$ws.Range('E41').Value = '  +4.24%  '

$ws.Range('E42').Value = '  +0.15%  '

$ws.Range('B43').Value = 'ONDO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D43' '1.20'
$ws.Range('E43').Value = '  +11.90%  '

$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D44' '4.39'
$ws.Range('E44').Value = '  +3.81%  '

$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D45' '41.29'
$ws.Range('E45').Value = '  -0.49%  '

$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D46' '0.758'
$ws.Range('E46').Value = '  -0.76%  '

Set-TextValue 'D47' '1.59'
$ws.Range('E47').Value = '  +3.64%  '

Set-TextValue 'D48' '23.00'
$ws.Range('E48').Value = '  +1.61%  '

Set-TextValue 'D49' '6.98'
$ws.Range('E49').Value = '  +5.77%  '

Set-TextValue 'D50' '23.18'
$ws.Range('E50').Value = '  +15.22%  '

Set-TextValue 'D51' '2.42'
$ws.Range('E51').Value = '  +13.70%  '
